$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll)
}

# 1. Update internship dates / add location
Replace-Text "June 2017 – present: Secret Escapes Software Development Internship" "June 2017 – September 2017: Secret Escapes Software Development Internship (London)"

# 2. Remove trailing period from IELTS line
Replace-Text "2015: IELTS: 7.5 (C1)." "2015: IELTS: 7.5 (C1)"

# 3. Delete the CISCO, Languages heading, and all language entries paragraphs
# (find "2014: CISCO IT Essentials." paragraph through "Japanese - Beginner." paragraph)
$ciscoPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "2014: CISCO IT Essentials.*") {
        $ciscoPara = $p
        break
    }
}
$japanesePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Japanese - Beginner.*") {
        $japanesePara = $p
        break
    }
}
$deleteRange = $d.Range($ciscoPara.Range.Start, $japanesePara.Range.End)
$deleteRange.Delete()

# 4. Update the Java professional skills line
Replace-Text "Java: multiple team projects written in it and IntelliJ, Gradle for building, Mockito and JUnit for mocking and unit testing and some other APIs where I used Git and GitHub / Bitbucket for VCS." "Java: multiple team projects written with it and IntelliJ, Gradle for building, Mockito and JUnit for mocking and unit testing and some other APIs where I used Git and GitHub / Bitbucket for VCS. Used it during my internship at Secret Escapes."

# 5. Replace the Scala/C++/HTML line with the Groovy/GRAILS line
Replace-Text "Scala, C++, HTML + CSS + Bootstrap + JavaScript + Angular: beginner, have done small projects with them." "Groovy and GRAILS: used these in a security-focused project in my first internship."

# 6. Replace the Software Engineering line with "Other languages / frameworks" line,
#    then insert a new paragraph after it with the original Software Engineering text + TDD sentence.
$sePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Software Engineering: familiar*") {
        $sePara = $p
        break
    }
}
$sePara.Range.Text = "Other languages / frameworks: Scala, Python, C++, JavaScript."
$sePara.Range.InsertParagraphAfter()
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Other languages / frameworks:*") {
        $newPara = $p
        break
    }
}
$afterPara = $newPara.Next()
$afterPara.Range.Text = "Software Engineering: familiar with Agile methodologies for project planning, such as Kanban and XP. I have applied my knowledge about software processes, architecture, and design to my projects. I did Test Driven Development during my internship at Secret Escapes."

# 7. Update the 2014 encryption algorithms line
Replace-Text "2014: I implemented several encryption algorithms, including: RSA, the Lorenz SZ42 cipher, and the MD5 hashing algorithm in C++ with the ttmath library for my college informatics attestation." "2014: I implemented several encryption algorithms, including: RSA, Enigma cipher, the Lorenz SZ42 cipher in C++ with the ttmath library for my high school Informatics attestation."

# 8. Update Sports line
Replace-Text "Sports: Swimming, Jogging, Muay Thai Kickboxing, weightlifting." "Sports: swimming, jogging, weightlifting."

# 9. Update Other hobbies line
Replace-Text "Other hobbies: reading, music, IT, gaming, chess, drawing." "Other hobbies: reading, music, tech, gaming."
